$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new translation rows 694-730 (37 rows x 3 columns: Key, Value, &Value)
$ws.Cells.Item(694, 1).Value = "statistics.totalFollowedUpFAdult"
$ws.Cells.Item(694, 2).Value = "Total Followed Up Female Adult Clients:"
$ws.Cells.Item(694, 3).Value = "&Total Followed Up Female Adult Clients:"
$ws.Cells.Item(695, 1).Value = "statistics.totalFollowedUpMAdult"
$ws.Cells.Item(695, 2).Value = "Total Followed Up Male Adult Clients:"
$ws.Cells.Item(695, 3).Value = "&Total Followed Up Male Adult Clients:"
$ws.Cells.Item(696, 1).Value = "statistics.totalFollowedUpFChild"
$ws.Cells.Item(696, 2).Value = "Total Followed Up Female Child Clients:"
$ws.Cells.Item(696, 3).Value = "&Total Followed Up Female Child Clients:"
$ws.Cells.Item(697, 1).Value = "statistics.totalFollowedUpMChild"
$ws.Cells.Item(697, 2).Value = "Total Followed Up Male Child Clients:"
$ws.Cells.Item(697, 3).Value = "&Total Followed Up Male Child Clients:"
$ws.Cells.Item(698, 1).Value = "statistics.followUpVisits"
$ws.Cells.Item(698, 2).Value = "Follow Up Visits:"
$ws.Cells.Item(698, 3).Value = "&Follow Up Visits:"
$ws.Cells.Item(699, 1).Value = "statistics.newClients"
$ws.Cells.Item(699, 2).Value = "New Clients:"
$ws.Cells.Item(699, 3).Value = "&New Clients:"
$ws.Cells.Item(700, 1).Value = "statistics.totalNewFAdult"
$ws.Cells.Item(700, 2).Value = "Total New Female Adult Clients:"
$ws.Cells.Item(700, 3).Value = "&Total New Female Adult Clients:"
$ws.Cells.Item(701, 1).Value = "statistics.totalNewMAdult"
$ws.Cells.Item(701, 2).Value = "Total New Male Adult Clients:"
$ws.Cells.Item(701, 3).Value = "&Total New Male Adult Clients:"
$ws.Cells.Item(702, 1).Value = "statistics.totalNewFChild"
$ws.Cells.Item(702, 2).Value = "Total New Female Child Clients:"
$ws.Cells.Item(702, 3).Value = "&Total New Female Child Clients:"
$ws.Cells.Item(703, 1).Value = "statistics.totalNewMChild"
$ws.Cells.Item(703, 2).Value = "Total New Male Child Clients:"
$ws.Cells.Item(703, 3).Value = "&Total New Male Child Clients:"
$ws.Cells.Item(704, 1).Value = "statistics.allChildren"
$ws.Cells.Item(704, 2).Value = "All Children "
$ws.Cells.Item(704, 3).Value = "&All Children "
$ws.Cells.Item(705, 1).Value = "statistics.allAdults"
$ws.Cells.Item(705, 2).Value = "All Adults"
$ws.Cells.Item(705, 3).Value = "&All Adults"
$ws.Cells.Item(706, 1).Value = "statistics.adult"
$ws.Cells.Item(706, 2).Value = "Adult"
$ws.Cells.Item(706, 3).Value = "&Adult"
$ws.Cells.Item(707, 1).Value = "statistics.child"
$ws.Cells.Item(707, 2).Value = "Child"
$ws.Cells.Item(707, 3).Value = "&Child"
$ws.Cells.Item(708, 1).Value = "statistics.age"
$ws.Cells.Item(708, 2).Value = "Age"
$ws.Cells.Item(708, 3).Value = "&Age"
$ws.Cells.Item(709, 1).Value = "statistics.filterByDemographic"
$ws.Cells.Item(709, 2).Value = "Filter by Demographic"
$ws.Cells.Item(709, 3).Value = "&Filter by Demographic "
$ws.Cells.Item(710, 1).Value = "statistics.femaleChild"
$ws.Cells.Item(710, 2).Value = "Female Child"
$ws.Cells.Item(710, 3).Value = "&Female Child"
$ws.Cells.Item(711, 1).Value = "statistics.maleChild"
$ws.Cells.Item(711, 2).Value = "Male Child"
$ws.Cells.Item(711, 3).Value = "&Male Child"
$ws.Cells.Item(712, 1).Value = "statistics.femaleAdult"
$ws.Cells.Item(712, 2).Value = "Female Adult"
$ws.Cells.Item(712, 3).Value = "&Female Adult"
$ws.Cells.Item(713, 1).Value = "statistics.maleAdult"
$ws.Cells.Item(713, 2).Value = "Male Adult"
$ws.Cells.Item(713, 3).Value = "&Male Adult"
$ws.Cells.Item(714, 1).Value = "statistics.totalFChild"
$ws.Cells.Item(714, 2).Value = "Total Female Children:"
$ws.Cells.Item(714, 3).Value = "&Total Female Children:"
$ws.Cells.Item(715, 1).Value = "statistics.totalMChild"
$ws.Cells.Item(715, 2).Value = "Total Male Children:"
$ws.Cells.Item(715, 3).Value = "&Total Male Children:"
$ws.Cells.Item(716, 1).Value = "statistics.totalFAdult"
$ws.Cells.Item(716, 2).Value = "Total Female Adults:"
$ws.Cells.Item(716, 3).Value = "&Total Female Adults:"
$ws.Cells.Item(717, 1).Value = "statistics.totalMAdult"
$ws.Cells.Item(717, 2).Value = "Total Male Adults:"
$ws.Cells.Item(717, 3).Value = "&Total Male Adults:"
$ws.Cells.Item(718, 1).Value = "statistics.totalFChildFollowUpVisits"
$ws.Cells.Item(718, 2).Value = "Total Female Children Follow Up Visits:"
$ws.Cells.Item(718, 3).Value = "&Total Female Children Follow Up Visits:"
$ws.Cells.Item(719, 1).Value = "statistics.totalMChildFollowUpVisits"
$ws.Cells.Item(719, 2).Value = "Total Male Children Follow Up Visits:"
$ws.Cells.Item(719, 3).Value = "&Total Male Children Follow Up Visits:"
$ws.Cells.Item(720, 1).Value = "statistics.totalFAdultFollowUpVisits"
$ws.Cells.Item(720, 2).Value = "Total Female Adult Follow Up Visits:"
$ws.Cells.Item(720, 3).Value = "&Total Female Adult Follow Up Visits:"
$ws.Cells.Item(721, 1).Value = "statistics.totalMAdultFollowUpVisits"
$ws.Cells.Item(721, 2).Value = "Total Male Adult Follow Up Visits:"
$ws.Cells.Item(721, 3).Value = "&Total Male Adult Follow Up Visits:"
$ws.Cells.Item(722, 1).Value = "statistics.selectAtLeastOne"
$ws.Cells.Item(722, 2).Value = "Select at least one Gender and Age option "
$ws.Cells.Item(722, 3).Value = "&Select at least one Gender and Age option "
$ws.Cells.Item(723, 1).Value = "statistics.warning"
$ws.Cells.Item(723, 2).Value = "Warning"
$ws.Cells.Item(723, 3).Value = "&Warning"
$ws.Cells.Item(724, 1).Value = "statistics.totalDisFChild"
$ws.Cells.Item(724, 2).Value = "Total Female Children With Disabilities: "
$ws.Cells.Item(724, 3).Value = "&Total Female Children With Disabilities: "
$ws.Cells.Item(725, 1).Value = "statistics.totalDisMChild"
$ws.Cells.Item(725, 2).Value = "Total Male Chidlren With Disabilities: "
$ws.Cells.Item(725, 3).Value = "&Total Male Chidlren With Disabilities: "
$ws.Cells.Item(726, 1).Value = "statistics.totalDisFAdult"
$ws.Cells.Item(726, 2).Value = "Total Female Adults With Disabilities: "
$ws.Cells.Item(726, 3).Value = "&Total Female Adults With Disabilities: "
$ws.Cells.Item(727, 1).Value = "statistics.totalDisMAdult"
$ws.Cells.Item(727, 2).Value = "Total Male Adults With Disabilities: "
$ws.Cells.Item(727, 3).Value = "&Total Male Adults With Disabilities: "
$ws.Cells.Item(728, 1).Value = "clientFields.hcrType"
$ws.Cells.Item(728, 2).Value = "HCR Type"
$ws.Cells.Item(728, 3).Value = "&HCR Type"
$ws.Cells.Item(729, 1).Value = "clientFields.hostCommunity"
$ws.Cells.Item(729, 2).Value = "Host Community"
$ws.Cells.Item(729, 3).Value = "&Host Community"
$ws.Cells.Item(730, 1).Value = "clientFields.refugee"
$ws.Cells.Item(730, 2).Value = "Refugee"
$ws.Cells.Item(730, 3).Value = "&Refugee"

# Apply the standard "translation row" formatting used elsewhere in the sheet
# (column A uses style from A649, columns B/C use style from B649:C649)
$ws.Range("A649").Copy()
$ws.Range("A694:A730").PasteSpecial(-4122)
$ws.Range("B649:C649").Copy()
$ws.Range("B694:C730").PasteSpecial(-4122)

# Row 731: a new, distinctly-styled empty marker cell (white fill, dark gray Arial 11 text, left aligned)
$marker = $ws.Cells.Item(731, 1)
$marker.Font.Name = "Arial"
$marker.Font.Size = 11
$marker.Font.Color = 2630431
$marker.Interior.Color = 16777215
$marker.HorizontalAlignment = -4131
